$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header row values for consistency (lowercase, no spaces)
$ws.Range("B1").Value = "firstname"
$ws.Range("C1").Value = "lastname"
$ws.Range("A1").Value = "patientid"

# Move selection to A2
$ws.Range("A2").Select()
